$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8..40 down to 9..41
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly price record
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = (Get-Date -Year 2021 -Month 10 -Day 12 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100112009
$ws.Range("G8").Value = "Acelga"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 950
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 975
$ws.Range("N8").Value = '$/atado 2,5 a 3 kilos'
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 325
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = "Hortaliza"

# Match the date cell formatting used by the other date cells in column D
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
